$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Unveiling the Enigma of Black Holes" "Exploring the Enigmatic Realm of Biochemistry: Unveiling the Secrets of Life"

# Body paragraph 1 (intro)
Replace-Text "Black holes, celestial enigmas that captivate the human imagination, stand as the ultimate cosmic paradox, embodying the interplay of immense density and profound emptiness" "Our world is a symphony of intricate processes, a tapestry woven with the threads of life"
Replace-Text " Their existence, predicted by the theories of general relativity, presents a challenging frontier, where the known laws of physics unravel and the very fabric of spacetime distorts" " In the heart of this symphony lies biochemistry, an enigmatic dance of molecules that unravels the very essence of existence"
Replace-Text " To fathom these cosmic behemoths and unmask the mysteries that shroud them, scientists embark on an intricate exploration, drawing on diverse fields of knowledge and experimental ingenuity" " From the smallest cellular structures to the vast ecosystems that shape our planet, biochemistry holds the key to understanding the profound mysteries of life"

Replace-Text "In the realm of astrophysics, observations conducted through telescopes, both ground-based and space-borne, offer glimpses into the enigmatic realm of black holes" "In this exploration, we embark on a journey through the enigmatic realm of biochemistry, unraveling the secrets of life's intricate processes"
Replace-Text " By meticulously analyzing the distortion of light and the emission of gravitational waves, astronomers piece together the characteristics and behavior of these celestial behemoths" " We delve into the structure and interactions of molecules, unveiling the symphony of chemical reactions that orchestrate the delicate balance of life"
Replace-Text " Their findings provide invaluable insights into the mechanics underlying the gravitational fields and accretion disks that surround black holes, unveiling the energetic processes that govern their extraordinary phenomena" " From the intricate dance of DNA replication to the intricate choreography of protein synthesis, we uncover the profound elegance of life's molecular mechanisms"

Replace-Text "Complementing observational efforts, theoretical physicists grapple with the complexities of mathematical models, seeking to unravel the enigmatic nature of black holes" "Furthermore, we venture into the enigmatic realm of metabolism, where the symphony of chemical reactions orchestrates the energy flow that powers life"
Replace-Text " Numerical simulations and analytical calculations probe the intricate interplay of gravity, quantum mechanics, and thermodynamics, venturing into unexplored territories of physics" " We unveil the intricate pathways of cellular respiration, photosynthesis, and glycolysis, uncovering the intricate choreography of molecules that provide the energy for life's myriad functions"
Replace-Text " These theoretical endeavors push the boundaries of human understanding, illuminating the fundamental properties and behaviors of black holes, and forging connections between disparate physical theories" " Through this journey, we unravel the profound interconnectedness of life, revealing the exquisite interplay between organisms and their environment"

# Summary paragraph
Replace-Text "Black holes, as intriguing as they are enigmatic, challenge our comprehension of the universe" "In the enchanting realm of biochemistry, we discover the symphony of molecular interactions that orchestrates the intricate dance of life"
Replace-Text " Through observations and theoretical explorations, scientists strive to unravel the mysteries surrounding these celestial phenomena" " From the structure of molecules to the choreography of chemical reactions, biochemistry unveils the profound elegance of life's mechanisms"
Replace-Text " The intricate interplay of gravity, quantum mechanics, and thermodynamics presents a profound scientific puzzle, captivating the imagination and driving the pursuit of knowledge" " We explore the enigmatic pathways of metabolism, revealing the energy flow that powers life and the exquisite interconnectedness of organisms and their environment"
